$wb = $excel.ActiveWorkbook

# --- Summary sheet: update aggregate stats after closing trade #36 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1198.34   # Current Capital
$wsSummary.Range("B4").Value = -1.66     # Total P&L $
$wsSummary.Range("B5").Value = -0.92     # Total P&L %
$wsSummary.Range("B6").Value = 36        # Total Trades
$wsSummary.Range("B8").Value = 19        # Losing Trades
$wsSummary.Range("B9").Value = 41.67     # Win Rate %

# --- Strategy Status sheet: update MarketMaking strategy row (row 4) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 98.34      # Capital
$wsStatus.Range("D4").Value = 36         # Trades
$wsStatus.Range("E4").Value = -1.66      # P&L $
$wsStatus.Range("F4").Value = -1.66      # P&L %
$wsStatus.Range("G4").Value = 41.67      # Win Rate %

# --- Append new closed trade (#36) to "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(37, 1).Value = 36
    $ws.Cells.Item(37, 2).Value = "'2026-02-17"
    $ws.Cells.Item(37, 3).Value = "'13:23:20"
    $ws.Cells.Item(37, 4).Value = "MarketMaking"
    $ws.Cells.Item(37, 5).Value = "UP"
    $ws.Cells.Item(37, 6).Value = 0.63
    $ws.Cells.Item(37, 7).Value = 0.03
    $ws.Cells.Item(37, 8).Value = "CLOSED"
    $ws.Cells.Item(37, 9).Value = -95.2381
    $ws.Cells.Item(37, 10).Value = -0.6
    $ws.Cells.Item(37, 11).Value = 98.34
    $ws.Cells.Item(37, 12).Value = 0
    $ws.Cells.Item(37, 13).Value = 0
    $ws.Cells.Item(37, 14).Value = 0.6
    $ws.Cells.Item(37, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(37, 16).Value = "early_exit"
    $ws.Cells.Item(37, 17).Value = 0.14
}
